# Update TPM values for Ntf3-Ntrk2 LR-pair data (rows 2-9) and remove the
# extra "MuSCs" sending-cluster block (rows 10-13), matching the refreshed
# NATMI TPM output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.798983
$ws.Range("H2").Value = 11.396949
$ws.Range("I2").Value = 0.9354058228349695
$ws.Range("J2").Value = 0.9354058228349694
$ws.Range("M2").Value = 0.6189250000000001
$ws.Range("N2").Value = 1.856775
$ws.Range("O2").Value = 0.09614699503454774
$ws.Range("P2").Value = 0.09614699503454775
$ws.Range("Q2").Value = 2.351285553275
$ws.Range("R2").Value = 21.161569979475
$ws.Range("S2").Value = 0.08993645900340086
$ws.Range("T2").Value = 0.08993645900340086
$ws.Range("G3").Value = 3.798983
$ws.Range("H3").Value = 11.396949
$ws.Range("I3").Value = 0.9354058228349695
$ws.Range("J3").Value = 0.9354058228349694
$ws.Range("O3").Value = 0.7879294335349575
$ws.Range("P3").Value = 0.7879294335349576
$ws.Range("Q3").Value = 19.26890271927066
$ws.Range("R3").Value = 173.420124473436
$ws.Range("S3").Value = 0.7370337801116583
$ws.Range("T3").Value = 0.7370337801116583
$ws.Range("G4").Value = 3.798983
$ws.Range("H4").Value = 11.396949
$ws.Range("I4").Value = 0.9354058228349695
$ws.Range("J4").Value = 0.9354058228349694
$ws.Range("M4").Value = 0.01220666666666667
$ws.Range("N4").Value = 0.03662
$ws.Range("O4").Value = 0.001896246426284896
$ws.Range("P4").Value = 0.001896246426284896
$ws.Range("Q4").Value = 0.04637291915333333
$ws.Range("R4").Value = 0.4173562723799999
$ws.Range("S4").Value = 0.001773759948676894
$ws.Range("T4").Value = 0.001773759948676894
$ws.Range("G5").Value = 3.798983
$ws.Range("H5").Value = 11.396949
$ws.Range("I5").Value = 0.9354058228349695
$ws.Range("J5").Value = 0.9354058228349694
$ws.Range("M5").Value = 0.7340256666666667
$ws.Range("N5").Value = 2.202077
$ws.Range("O5").Value = 0.1140273250042099
$ws.Range("P5").Value = 0.1140273250042099
$ws.Range("Q5").Value = 2.788551029230333
$ws.Range("R5").Value = 25.096959263073
$ws.Range("S5").Value = 0.1066618237712334
$ws.Range("T5").Value = 0.1066618237712334
$ws.Range("I6").Value = 0.06459417716503056
$ws.Range("J6").Value = 0.06459417716503056
$ws.Range("M6").Value = 0.6189250000000001
$ws.Range("N6").Value = 1.856775
$ws.Range("O6").Value = 0.09614699503454774
$ws.Range("P6").Value = 0.09614699503454775
$ws.Range("Q6").Value = 0.1623673403416667
$ws.Range("R6").Value = 1.461306063075
$ws.Range("S6").Value = 0.006210536031146891
$ws.Range("T6").Value = 0.006210536031146891
$ws.Range("I7").Value = 0.06459417716503056
$ws.Range("J7").Value = 0.06459417716503056
$ws.Range("O7").Value = 0.7879294335349575
$ws.Range("P7").Value = 0.7879294335349576
$ws.Range("S7").Value = 0.05089565342329921
$ws.Range("T7").Value = 0.05089565342329922
$ws.Range("I8").Value = 0.06459417716503056
$ws.Range("J8").Value = 0.06459417716503056
$ws.Range("M8").Value = 0.01220666666666667
$ws.Range("N8").Value = 0.03662
$ws.Range("O8").Value = 0.001896246426284896
$ws.Range("P8").Value = 0.001896246426284896
$ws.Range("Q8").Value = 0.003202268451111111
$ws.Range("R8").Value = 0.02882041606
$ws.Range("S8").Value = 0.0001224864776080026
$ws.Range("T8").Value = 0.0001224864776080027
$ws.Range("I9").Value = 0.06459417716503056
$ws.Range("J9").Value = 0.06459417716503056
$ws.Range("M9").Value = 0.7340256666666667
$ws.Range("N9").Value = 2.202077
$ws.Range("O9").Value = 0.1140273250042099
$ws.Range("P9").Value = 0.1140273250042099
$ws.Range("Q9").Value = 0.1925625806667778
$ws.Range("R9").Value = 1.733063226001
$ws.Range("S9").Value = 0.007365501232976451
$ws.Range("T9").Value = 0.007365501232976452

# Remove the trailing MuSCs sending-cluster rows (10-13); data now ends at row 9.
$ws.Range("A10:T13").Delete()
